$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.574307304785894
$ws.Range("C2").Value = 0.0163727959697733
$ws.Range("D2").Value = 0.0604534005037783
$ws.Range("E2").Value = 0.945843828715365
$ws.Range("F2").Value = 0.0113350125944584
$ws.Range("G2").Value = 0.944584382871537
$ws.Range("H2").Value = 0.0138539042821159
$ws.Range("I2").Value = 0.730478589420655
$ws.Range("J2").Value = 0.0428211586901763
$ws.Range("K2").Value = 0.0390428211586902
$ws.Range("L2").Value = 0.0377833753148615
$ws.Range("M2").Value = 0.838790931989924
$ws.Range("N2").Value = 0.0100755667506297
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.0327455919395466
$ws.Range("Q2").Value = 0.928211586901763
$ws.Range("R2").Value = 0.00629722921914358
$ws.Range("S2").Value = 0.00377833753148615
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0.00629722921914358
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0.0113350125944584
$ws.Range("X2").Value = 0.00755667506297229
$ws.Range("B3").Value = 0.0654911838790932
$ws.Range("C3").Value = 0.346347607052897
$ws.Range("D3").Value = 0.840050377833753
$ws.Range("E3").Value = 0.044080604534005
$ws.Range("F3").Value = 0.00125944584382872
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.0705289672544081
$ws.Range("I3").Value = 0.0377833753148615
$ws.Range("J3").Value = 0.132241813602015
$ws.Range("K3").Value = 0.919395465994962
$ws.Range("L3").Value = 0.958438287153652
$ws.Range("M3").Value = 0.109571788413098
$ws.Range("N3").Value = 0.919395465994962
$ws.Range("O3").Value = 0.0151133501259446
$ws.Range("P3").Value = 0.00125944584382872
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.851385390428212
$ws.Range("S3").Value = 0.983627204030227
$ws.Range("T3").Value = 0.0554156171284635
$ws.Range("U3").Value = 0.00755667506297229
$ws.Range("V3").Value = 0.0214105793450882
$ws.Range("W3").Value = 0.0327455919395466
$ws.Range("X3").Value = 0.00125944584382872
$ws.Range("B4").Value = 0.353904282115869
$ws.Range("C4").Value = 0.0251889168765743
$ws.Range("D4").Value = 0.00377833753148615
$ws.Range("E4").Value = 0.00377833753148615
$ws.Range("F4").Value = 0.945843828715365
$ws.Range("G4").Value = 0.0516372795969773
$ws.Range("H4").Value = 0.00251889168765743
$ws.Range("I4").Value = 0.0188916876574307
$ws.Range("J4").Value = 0.052896725440806
$ws.Range("K4").Value = 0.0390428211586902
$ws.Range("L4").Value = 0.00125944584382872
$ws.Range("M4").Value = 0.00629722921914358
$ws.Range("N4").Value = 0.00377833753148615
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.963476070528967
$ws.Range("Q4").Value = 0.00377833753148615
$ws.Range("R4").Value = 0.134760705289673
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0.00251889168765743
$ws.Range("V4").Value = 0.0541561712846348
$ws.Range("W4").Value = 0.953400503778338
$ws.Range("X4").Value = 0.954659949622166
$ws.Range("B5").Value = 0.00629722921914358
$ws.Range("C5").Value = 0.612090680100756
$ws.Range("D5").Value = 0.0931989924433249
$ws.Range("E5").Value = 0.00503778337531486
$ws.Range("F5").Value = 0.0415617128463476
$ws.Range("G5").Value = 0.00377833753148615
$ws.Range("H5").Value = 0.913098236775819
$ws.Range("I5").Value = 0.212846347607053
$ws.Range("J5").Value = 0.772040302267002
$ws.Range("K5").Value = 0.00125944584382872
$ws.Range("L5").Value = 0.00251889168765743
$ws.Range("M5").Value = 0.0453400503778338
$ws.Range("N5").Value = 0.0667506297229219
$ws.Range("O5").Value = 0.984886649874055
$ws.Range("P5").Value = 0.00251889168765743
$ws.Range("Q5").Value = 0.0680100755667506
$ws.Range("R5").Value = 0.00755667506297229
$ws.Range("S5").Value = 0.0125944584382872
$ws.Range("T5").Value = 0.944584382871537
$ws.Range("U5").Value = 0.983627204030227
$ws.Range("V5").Value = 0.924433249370277
$ws.Range("W5").Value = 0.00251889168765743
$ws.Range("X5").Value = 0.0365239294710327